$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ENTER INTO UI" input row (row 2) with the new resource ---
$ws.Range("A2").Value = "AOD_NUAspanish"
$ws.Range("B2").Value = "Never Use Alone (Spanish)"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").Value = "Never Use Alone (Spanish)"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").Value = "T"
$ws.Range("F2").Value = "T"

# --- Remove the old hard-coded accordion/server/info rows (20:33) ---
$ws.Range("A20:A33").EntireRow.Delete()

# --- Populate the new data-entry row (row 13) (C13 is handled separately: it
# --- becomes a hyperlink cell, which carries its own built-in style) ---
$row13a = $ws.Range("A13:B13")
$row13a.NumberFormat = "@"
$row13a.Font.Name = "Times New Roman"
$row13a.Font.Size = 12
$row13a.Font.Bold = $false

$row13b = $ws.Range("D13:I13")
$row13b.NumberFormat = "@"
$row13b.Font.Name = "Times New Roman"
$row13b.Font.Size = 12
$row13b.Font.Bold = $false

$ws.Range("A13").Value = "Never Use Alone (Spanish)"
$ws.Range("C13").Value = "tel:18009285330"
$ws.Range("D13").Value = "https://neverusealone.com/"
$ws.Range("G13").Value = "AOD"
$ws.Range("H13").Value = "Off"
$ws.Range("I13").Value = "National anonymous hotline - Will stay on the phone with you if you use alone and will call for help to your address if you experience an overdose - Overdose prevention"

# --- Hyperlink the phone number cell (applies the built-in Hyperlink style) ---
$ws.Hyperlinks.Add($ws.Range("C13"), "tel:18009285330")

# --- Populate the generated UI snippets (rows 16:18) ---
$ws.Range("A16").Value = "mod_Accordion_ui('AOD_NUAspanish')"
$ws.Range("A17").Value = "mod_Accordion_server('AOD_NUAspanish', selector=selection, data=AODdata, title = c('Never Use Alone (Spanish)'), Visible = T)"
$ws.Range("A18").Value = "mod_info_server('AOD_NUAspanish', selector = selection, data = AODdata, rownametitle = c('Never Use Alone (Spanish)'), phone = T, website = T)"

$wb.Save()
